$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2033898305084746
$ws.Range("C2").Value = 0.5466101694915254
$ws.Range("J2").Value = 0.03813559322033899
$ws.Range("O2").Value = 0.00423728813559322
$ws.Range("P2").Value = 0.1059322033898305
$ws.Range("S2").Value = 0.1016949152542373
$ws.Range("B3").Value = 0.006993006993006993
$ws.Range("C3").Value = 0.04895104895104895
$ws.Range("J3").Value = 0.04195804195804196
$ws.Range("P3").Value = 0.7552447552447552
$ws.Range("S3").Value = 0.1468531468531468
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.71875
$ws.Range("S4").Value = 0.21875
$ws.Range("B6").Value = 0.02415458937198068
$ws.Range("E6").Value = 0.004830917874396135
$ws.Range("F6").Value = 0.07729468599033816
$ws.Range("J6").Value = 0.2801932367149759
$ws.Range("O6").Value = 0.02415458937198068
$ws.Range("Q6").Value = 0.1835748792270532
$ws.Range("R6").Value = 0.04347826086956522
$ws.Range("S6").Value = 0.3623188405797101
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("D7").Value = 0.01449275362318841
$ws.Range("E7").Value = 0.004830917874396135
$ws.Range("F7").Value = 0.04347826086956522
$ws.Range("J7").Value = 0.178743961352657
$ws.Range("O7").Value = 0.01449275362318841
$ws.Range("Q7").Value = 0.1739130434782609
$ws.Range("R7").Value = 0.0821256038647343
$ws.Range("S7").Value = 0.3768115942028986
$ws.Range("B8").Value = 0.05853658536585366
$ws.Range("D8").Value = 0.01707317073170732
$ws.Range("F8").Value = 0.06585365853658537
$ws.Range("J8").Value = 0.1463414634146341
$ws.Range("O8").Value = 0.02195121951219512
$ws.Range("Q8").Value = 0.1780487804878049
$ws.Range("R8").Value = 0.1317073170731707
$ws.Range("S8").Value = 0.3804878048780488
$ws.Range("B9").Value = 0.08695652173913043
$ws.Range("D9").Value = 0.01739130434782609
$ws.Range("E9").Value = 0.004347826086956522
$ws.Range("F9").Value = 0.03043478260869565
$ws.Range("J9").Value = 0.1260869565217391
$ws.Range("O9").Value = 0.01739130434782609
$ws.Range("Q9").Value = 0.1695652173913043
$ws.Range("R9").Value = 0.08260869565217391
$ws.Range("S9").Value = 0.4652173913043478
$ws.Range("B10").Value = 0.0936205468102734
$ws.Range("D10").Value = 0.01491300745650373
$ws.Range("E10").Value = 0.001657000828500414
$ws.Range("F10").Value = 0.06710853355426678
$ws.Range("J10").Value = 0.1425020712510356
$ws.Range("O10").Value = 0.01077050538525269
$ws.Range("Q10").Value = 0.2162386081193041
$ws.Range("R10").Value = 0.07787903893951947
$ws.Range("S10").Value = 0.3753106876553438
$ws.Range("G11").Value = 0.1107142857142857
$ws.Range("J11").Value = 0.06785714285714285
$ws.Range("K11").Value = 0.1678571428571428
$ws.Range("L11").Value = 0.6214285714285714
$ws.Range("S11").Value = 0.03214285714285714
$ws.Range("G12").Value = 0.8108108108108109
$ws.Range("J12").Value = 0.1243243243243243
$ws.Range("K12").Value = 0.01081081081081081
$ws.Range("L12").Value = 0.03243243243243243
$ws.Range("S12").Value = 0.02162162162162162
$ws.Range("G13").Value = 0.7441860465116279
$ws.Range("J13").Value = 0.2093023255813954
$ws.Range("S13").Value = 0.04651162790697674
$ws.Range("F15").Value = 0.04366812227074236
$ws.Range("H15").Value = 0.1135371179039301
$ws.Range("I15").Value = 0.09606986899563319
$ws.Range("J15").Value = 0.3187772925764192
$ws.Range("K15").Value = 0.09170305676855896
$ws.Range("M15").Value = 0.01310043668122271
$ws.Range("O15").Value = 0.09606986899563319
$ws.Range("S15").Value = 0.2270742358078603
$ws.Range("F16").Value = 0.006622516556291391
$ws.Range("H16").Value = 0.2119205298013245
$ws.Range("I16").Value = 0.09933774834437085
$ws.Range("J16").Value = 0.3774834437086093
$ws.Range("K16").Value = 0.09933774834437085
$ws.Range("M16").Value = 0.03973509933774835
$ws.Range("N16").Value = 0.006622516556291391
$ws.Range("O16").Value = 0.05960264900662252
$ws.Range("S16").Value = 0.09933774834437085
$ws.Range("F17").Value = 0.02040816326530612
$ws.Range("H17").Value = 0.1791383219954649
$ws.Range("I17").Value = 0.1179138321995465
$ws.Range("J17").Value = 0.3718820861678004
$ws.Range("K17").Value = 0.09523809523809523
$ws.Range("M17").Value = 0.02947845804988662
$ws.Range("N17").Value = 0.002267573696145125
$ws.Range("O17").Value = 0.07482993197278912
$ws.Range("S17").Value = 0.108843537414966
$ws.Range("F18").Value = 0.02631578947368421
$ws.Range("H18").Value = 0.1947368421052632
$ws.Range("I18").Value = 0.131578947368421
$ws.Range("J18").Value = 0.3842105263157894
$ws.Range("K18").Value = 0.08421052631578947
$ws.Range("M18").Value = 0.01578947368421053
$ws.Range("O18").Value = 0.04736842105263158
$ws.Range("S18").Value = 0.1157894736842105
$ws.Range("F19").Value = 0.01094276094276094
$ws.Range("H19").Value = 0.1944444444444444
$ws.Range("I19").Value = 0.09764309764309764
$ws.Range("J19").Value = 0.3779461279461279
$ws.Range("K19").Value = 0.111952861952862
$ws.Range("M19").Value = 0.01683501683501683
$ws.Range("N19").Value = 0.0008417508417508417
$ws.Range("O19").Value = 0.08080808080808081
$ws.Range("S19").Value = 0.1085858585858586
